# Applies the "provisional" data update to the delivery route sheet:
# - Column E width changes from 17 to 13
# - Rows 2-16: Nº Pedido (A), Data de entrega (B), Bairro (E) are updated
#   for every row; row 6's Período (C) changes from "manhã" to "tarde"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column E (Bairro) from width 17 to 13
$ws.Columns.Item(5).ColumnWidth = 13

$data = @(
    @{ Row = 2;  A = "11185"; B = "19/11/2024"; E = "Passa Vinte" },
    @{ Row = 3;  A = "11185"; B = "19/11/2024"; E = "Passa Vinte" },
    @{ Row = 4;  A = "11075"; B = "19/11/2024"; E = "Bela Vista" },
    @{ Row = 5;  A = "11075"; B = "19/11/2024"; E = "Bela Vista" },
    @{ Row = 6;  A = "11073"; B = "19/11/2024"; C = "tarde"; E = "Bela Vista" },
    @{ Row = 7;  A = "11073"; B = "19/11/2024"; E = "Bela Vista" },
    @{ Row = 8;  A = "11204"; B = "19/11/2024"; E = "Potecas" },
    @{ Row = 9;  A = "11204"; B = "19/11/2024"; E = "Potecas" },
    @{ Row = 10; A = "11083"; B = "19/11/2024"; E = "Brejarú" },
    @{ Row = 11; A = "11057"; B = "21/11/2024"; E = "Aririú" },
    @{ Row = 12; A = "11057"; B = "21/11/2024"; E = "Aririú" },
    @{ Row = 13; A = "11214"; B = "21/11/2024"; E = "Rio Grande" },
    @{ Row = 14; A = "11181"; B = "21/11/2024"; E = "Pachecos" },
    @{ Row = 15; A = "11181"; B = "21/11/2024"; E = "Pachecos" },
    @{ Row = 16; A = "11159"; B = "21/11/2024"; E = "Furadinho" }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    if ($item.ContainsKey("C")) {
        $ws.Cells.Item($r, 3).Value = $item.C
    }
    $ws.Cells.Item($r, 5).Value = $item.E
}
